# Auto-generated Excel COM-interop script applying value updates
# to the Sheets workbook, per the provided unified diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1067
$ws.Range("I2").Value = 270.8
$ws.Range("K2").Value = 270.8
$ws.Range("M2").Value = -157.8
$ws.Range("H4").Value = 214.22223
$ws.Range("I4").Value = 214.22223
$ws.Range("K4").Value = 214.22223
$ws.Range("M4").Value = -100.22223
$ws.Range("H6").Value = 224.36842
$ws.Range("I6").Value = 80
$ws.Range("J6").Value = 275.92856
$ws.Range("K6").Value = 240
$ws.Range("L6").Value = 827.78568
$ws.Range("M6").Value = -128
$ws.Range("N6").Value = -1051.78568
$ws.Range("H12").Value = 608.35
$ws.Range("I12").Value = 316.8125
$ws.Range("J12").Value = 1774.5
$ws.Range("K12").Value = 316.8125
$ws.Range("L12").Value = 1774.5
$ws.Range("M12").Value = -146.8125
$ws.Range("N12").Value = -2114.5
$ws.Range("H32").Value = 3997.6667
$ws.Range("J32").Value = 4496.5
$ws.Range("L32").Value = 4496.5
$ws.Range("N32").Value = -5148.5
$ws.Range("H33").Value = 1178.8096
$ws.Range("I33").Value = 1338.2222
$ws.Range("J33").Value = 222.33333
$ws.Range("K33").Value = 1338.2222
$ws.Range("L33").Value = 222.33333
$ws.Range("M33").Value = -1109.2222
$ws.Range("N33").Value = -680.3333299999999
$ws.Range("H41").Value = 614.5714
$ws.Range("I41").Value = 466.36365
$ws.Range("J41").Value = 777.6
$ws.Range("K41").Value = 466.36365
$ws.Range("L41").Value = 777.6
$ws.Range("M41").Value = -26.36365000000001
$ws.Range("N41").Value = -1657.6
$ws.Range("H42").Value = 74.27273
$ws.Range("I42").Value = 74.111115
$ws.Range("K42").Value = 222.333345
$ws.Range("M42").Value = 7.666654999999992
$ws.Range("H43").Value = 9786.6
$ws.Range("J43").Value = 8483.5
$ws.Range("L43").Value = 8483.5
$ws.Range("N43").Value = -8621.5
$ws.Range("H53").Value = 995.875
$ws.Range("I53").Value = 114
$ws.Range("J53").Value = 1877.75
$ws.Range("K53").Value = 114
$ws.Range("L53").Value = 1877.75
$ws.Range("M53").Value = 523
$ws.Range("N53").Value = -3151.75
$ws.Range("H58").Value = 174.61539
$ws.Range("J58").Value = 64
$ws.Range("L58").Value = 192
$ws.Range("N58").Value = -492
$ws.Range("H69").Value = 20258
$ws.Range("I69").Value = 20911
$ws.Range("K69").Value = 62733
$ws.Range("M69").Value = -61859
$ws.Range("H72").Value = 20258
$ws.Range("I72").Value = 20911
$ws.Range("K72").Value = 188199
$ws.Range("M72").Value = -183831
$ws.Range("H74").Value = 8395.091
$ws.Range("I74").Value = 9406.714
$ws.Range("K74").Value = 9406.714
$ws.Range("M74").Value = -8470.714
$ws.Range("H77").Value = 8395.091
$ws.Range("I77").Value = 9406.714
$ws.Range("K77").Value = 47033.57
$ws.Range("M77").Value = -42353.57
$ws.Range("H93").Value = 22100
$ws.Range("J93").Value = 22100
$ws.Range("L93").Value = 22100
$ws.Range("N93").Value = -27092
$ws.Range("H100").Value = 8468.75
$ws.Range("I100").Value = 1848.375
$ws.Range("K100").Value = 1848.375
$ws.Range("M100").Value = -1307.375
$ws.Range("H101").Value = 1830.1666
$ws.Range("I101").Value = 1183.4286
$ws.Range("J101").Value = 2735.6
$ws.Range("K101").Value = 3550.2858
$ws.Range("L101").Value = 8206.799999999999
$ws.Range("M101").Value = -1928.2858
$ws.Range("N101").Value = -11450.8
$ws.Range("H113").Value = 5342.1514
$ws.Range("I113").Value = 4805.85
$ws.Range("J113").Value = 6167.231
$ws.Range("K113").Value = 4805.85
$ws.Range("L113").Value = 6167.231
$ws.Range("M113").Value = -1551.85
$ws.Range("N113").Value = -12675.231
$ws.Range("H125").Value = 1220.6
$ws.Range("I125").Value = 1010.3333
$ws.Range("K125").Value = 9092.9997
$ws.Range("M125").Value = -6632.9997
$ws.Range("H130").Value = 100000
$ws.Range("J130").Value = 100000
$ws.Range("L130").Value = 100000
$ws.Range("N130").Value = -110040
$ws.Range("H132").Value = 11631087
$ws.Range("I132").Value = 14708858
$ws.Range("J132").Value = 3949.9443
$ws.Range("K132").Value = 44126574
$ws.Range("L132").Value = 11849.8329
$ws.Range("M132").Value = -44124044
$ws.Range("N132").Value = -16909.8329
$ws.Range("H135").Value = 10452.842
$ws.Range("I135").Value = 3152
$ws.Range("K135").Value = 28368
$ws.Range("M135").Value = -25833
$ws.Range("H137").Value = 1718.3077
$ws.Range("I137").Value = 1133.64
$ws.Range("K137").Value = 3400.92
$ws.Range("M137").Value = -850.9200000000001
$ws.Range("H138").Value = 4780.6
$ws.Range("I138").Value = 1406.5714
$ws.Range("J138").Value = 7732.875
$ws.Range("K138").Value = 4219.7142
$ws.Range("L138").Value = 23198.625
$ws.Range("M138").Value = 920.2857999999997
$ws.Range("N138").Value = -33478.625
$ws.Range("H141").Value = 3063.3928
$ws.Range("I141").Value = 2643.1
$ws.Range("J141").Value = 4114.125
$ws.Range("K141").Value = 7929.299999999999
$ws.Range("L141").Value = 12342.375
$ws.Range("M141").Value = -2749.299999999999
$ws.Range("N141").Value = -22702.375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1627.9117
$ws.Range("I2").Value = 1407.0952
$ws.Range("J2").Value = 1984.6154
$ws.Range("K2").Value = 1407.0952
$ws.Range("L2").Value = 1984.6154
$ws.Range("M2").Value = -1294.0952
$ws.Range("N2").Value = -2210.6154
$ws.Range("H32").Value = 2531.0156
$ws.Range("I32").Value = 2618.1064
$ws.Range("J32").Value = 2290.2354
$ws.Range("K32").Value = 2618.1064
$ws.Range("L32").Value = 2290.2354
$ws.Range("M32").Value = -2331.1064
$ws.Range("N32").Value = -2864.2354
$ws.Range("H61").Value = 5770.759
$ws.Range("I61").Value = 2763.3704
$ws.Range("K61").Value = 2763.3704
$ws.Range("M61").Value = -2551.3704
$ws.Range("H63").Value = 4728.4
$ws.Range("I63").Value = 4728.4
$ws.Range("K63").Value = 4728.4
$ws.Range("M63").Value = -4042.4
$ws.Range("H66").Value = 4728.4
$ws.Range("I66").Value = 4728.4
$ws.Range("K66").Value = 23642
$ws.Range("M66").Value = -20210
$ws.Range("H102").Value = 18182858
$ws.Range("I102").Value = 1219.8572
$ws.Range("K102").Value = 1219.8572
$ws.Range("M102").Value = 402.1428000000001
$ws.Range("H110").Value = 1609
$ws.Range("I110").Value = 1261.9546
$ws.Range("J110").Value = 2881.5
$ws.Range("K110").Value = 1261.9546
$ws.Range("L110").Value = 2881.5
$ws.Range("M110").Value = 783.0454
$ws.Range("N110").Value = -6971.5
$ws.Range("H116").Value = 1627.9117
$ws.Range("I116").Value = 1407.0952
$ws.Range("J116").Value = 1984.6154
$ws.Range("K116").Value = 1407.0952
$ws.Range("L116").Value = 1984.6154
$ws.Range("M116").Value = 886.9048
$ws.Range("N116").Value = -6572.6154
$ws.Range("H132").Value = 7111.1113
$ws.Range("I132").Value = 3288.7
$ws.Range("K132").Value = 9866.099999999999
$ws.Range("M132").Value = -7336.099999999999
$ws.Range("H136").Value = 5770.759
$ws.Range("I136").Value = 2763.3704
$ws.Range("K136").Value = 8290.111199999999
$ws.Range("M136").Value = -5740.111199999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1627.9117
$ws.Range("I3").Value = 1407.0952
$ws.Range("J3").Value = 1984.6154
$ws.Range("K3").Value = 1407.0952
$ws.Range("L3").Value = 1984.6154
$ws.Range("M3").Value = -1293.0952
$ws.Range("N3").Value = -2212.6154
$ws.Range("H5").Value = 1312.375
$ws.Range("I5").Value = 1166.6666
$ws.Range("J5").Value = 1749.5
$ws.Range("K5").Value = 1166.6666
$ws.Range("L5").Value = 1749.5
$ws.Range("M5").Value = -1053.6666
$ws.Range("N5").Value = -1975.5
$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -32122
$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -100608
$ws.Range("H88").Value = 21855.8
$ws.Range("J88").Value = 21855.8
$ws.Range("L88").Value = 21855.8
$ws.Range("N88").Value = -22667.8
$ws.Range("H91").Value = 21855.8
$ws.Range("J91").Value = 21855.8
$ws.Range("L91").Value = 21855.8
$ws.Range("N91").Value = -24663.8
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H99").Value = 2597
$ws.Range("I99").Value = 2597
$ws.Range("K99").Value = 2597
$ws.Range("M99").Value = -1099
$ws.Range("H131").Value = 21500
$ws.Range("J131").Value = 21500
$ws.Range("L131").Value = 21500
$ws.Range("N131").Value = -31580
$ws.Range("H132").Value = 87000
$ws.Range("J132").Value = 87000
$ws.Range("L132").Value = 87000
$ws.Range("N132").Value = -97120
$ws.Range("H134").Value = 2231.7097
$ws.Range("I134").Value = 1893.3
$ws.Range("J134").Value = 2847
$ws.Range("K134").Value = 5679.9
$ws.Range("L134").Value = 8541
$ws.Range("M134").Value = -3144.9
$ws.Range("N134").Value = -13611

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 252504500
$ws.Range("J4").Value = 336669340
$ws.Range("L4").Value = 336669340
$ws.Range("N4").Value = -336669564
$ws.Range("H22").Value = 600
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 3330.8918
$ws.Range("I31").Value = 2148.875
$ws.Range("J31").Value = 4231.476
$ws.Range("K31").Value = 2148.875
$ws.Range("L31").Value = 4231.476
$ws.Range("M31").Value = -1853.875
$ws.Range("N31").Value = -4821.476
$ws.Range("H34").Value = 3330.8918
$ws.Range("I34").Value = 2148.875
$ws.Range("J34").Value = 4231.476
$ws.Range("K34").Value = 2148.875
$ws.Range("L34").Value = 4231.476
$ws.Range("M34").Value = -1946.875
$ws.Range("N34").Value = -4635.476
$ws.Range("H58").Value = 3011
$ws.Range("I58").Value = 3532.1667
$ws.Range("J58").Value = 2663.5557
$ws.Range("K58").Value = 3532.1667
$ws.Range("L58").Value = 2663.5557
$ws.Range("M58").Value = -3329.1667
$ws.Range("N58").Value = -3069.5557
$ws.Range("H86").Value = 13089.2
$ws.Range("I86").Value = 9631
$ws.Range("J86").Value = 14571.286
$ws.Range("K86").Value = 9631
$ws.Range("L86").Value = 14571.286
$ws.Range("M86").Value = -8508
$ws.Range("N86").Value = -16817.286
$ws.Range("H89").Value = 13089.2
$ws.Range("I89").Value = 9631
$ws.Range("J89").Value = 14571.286
$ws.Range("K89").Value = 48155
$ws.Range("L89").Value = 72856.42999999999
$ws.Range("M89").Value = -42539
$ws.Range("N89").Value = -84088.42999999999
$ws.Range("H99").Value = 4581
$ws.Range("I99").Value = 4581
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4581
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3083
$ws.Range("N99").ClearContents()
$ws.Range("H122").Value = 5994.2144
$ws.Range("I122").Value = 3624
$ws.Range("J122").Value = 9154.5
$ws.Range("K122").Value = 10872
$ws.Range("L122").Value = 27463.5
$ws.Range("M122").Value = -8422
$ws.Range("N122").Value = -32363.5
$ws.Range("H126").Value = 4581
$ws.Range("I126").Value = 4581
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 13743
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11273
$ws.Range("N126").ClearContents()
$ws.Range("H130").Value = 119067.75
$ws.Range("J130").Value = 119067.75
$ws.Range("L130").Value = 119067.75
$ws.Range("N130").Value = -129107.75
$ws.Range("H132").Value = 3590.2563
$ws.Range("I132").Value = 2167.3667
$ws.Range("K132").Value = 6502.1001
$ws.Range("M132").Value = -3972.1001
$ws.Range("H134").Value = 1634.4
$ws.Range("I134").Value = 1563.1765
$ws.Range("K134").Value = 4689.529500000001
$ws.Range("M134").Value = -2154.529500000001
$ws.Range("H136").Value = 3011
$ws.Range("I136").Value = 3532.1667
$ws.Range("J136").Value = 2663.5557
$ws.Range("K136").Value = 10596.5001
$ws.Range("L136").Value = 7990.6671
$ws.Range("M136").Value = -8046.500100000001
$ws.Range("N136").Value = -13090.6671
$ws.Range("H138").Value = 87000
$ws.Range("J138").Value = 87000
$ws.Range("L138").Value = 87000
$ws.Range("N138").Value = -97280
$ws.Range("H141").Value = 133992.67
$ws.Range("J141").Value = 133992.67
$ws.Range("L141").Value = 133992.67
$ws.Range("N141").Value = -144352.67

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1552.4286
$ws.Range("I3").Value = 1552.4286
$ws.Range("K3").Value = 4657.2858
$ws.Range("M3").Value = -4545.2858
$ws.Range("H7").Value = 94.888885
$ws.Range("I7").Value = 94.888885
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 284.666655
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -172.666655
$ws.Range("N7").ClearContents()
$ws.Range("H34").Value = 1493
$ws.Range("J34").Value = 1183.5
$ws.Range("L34").Value = 3550.5
$ws.Range("N34").Value = -3718.5
$ws.Range("H39").Value = 4199
$ws.Range("J39").Value = 4954.615
$ws.Range("L39").Value = 14863.845
$ws.Range("N39").Value = -15451.845
$ws.Range("H40").Value = 118.53846
$ws.Range("I40").Value = 126.333336
$ws.Range("J40").Value = 25
$ws.Range("K40").Value = 505.333344
$ws.Range("L40").Value = 100
$ws.Range("M40").Value = -436.333344
$ws.Range("N40").Value = -238
$ws.Range("H46").Value = 2304
$ws.Range("J46").Value = 350
$ws.Range("L46").Value = 1050
$ws.Range("N46").Value = -1232
$ws.Range("H55").Value = 6765.5625
$ws.Range("I55").Value = 1499.5
$ws.Range("J55").Value = 7517.857
$ws.Range("K55").Value = 4498.5
$ws.Range("L55").Value = 22553.571
$ws.Range("M55").Value = -4321.5
$ws.Range("N55").Value = -22907.571
$ws.Range("H60").Value = 758.1429000000001
$ws.Range("I60").Value = 151.75
$ws.Range("J60").Value = 1566.6666
$ws.Range("K60").Value = 455.25
$ws.Range("L60").Value = 4699.9998
$ws.Range("M60").Value = -204.25
$ws.Range("N60").Value = -5201.9998
$ws.Range("H114").Value = 1555.7142
$ws.Range("J114").Value = 1694.6666
$ws.Range("L114").Value = 5083.9998
$ws.Range("N114").Value = -11591.9998
$ws.Range("H122").Value = 937
$ws.Range("J122").Value = 966.75
$ws.Range("L122").Value = 8700.75
$ws.Range("N122").Value = -13600.75
$ws.Range("H139").Value = 3218.0476
$ws.Range("I139").Value = 2438.6
$ws.Range("K139").Value = 7315.799999999999
$ws.Range("M139").Value = -2175.799999999999
$ws.Range("H140").Value = 4142.727
$ws.Range("I140").Value = 3162.1538
$ws.Range("K140").Value = 9486.4614
$ws.Range("M140").Value = -4306.4614

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 50225.816
$ws.Range("J95").Value = 50225.816
$ws.Range("L95").Value = 50225.816
$ws.Range("N95").Value = -55717.816
$ws.Range("H101").Value = 52228
$ws.Range("J101").Value = 52228
$ws.Range("L101").Value = 52228
$ws.Range("N101").Value = -58718
$ws.Range("H102").Value = 22875.625
$ws.Range("I102").Value = 2155.318
$ws.Range("K102").Value = 2155.318
$ws.Range("M102").Value = -533.3180000000002
$ws.Range("H113").Value = 3096.5833
$ws.Range("I113").Value = 3312
$ws.Range("J113").Value = 2665.75
$ws.Range("K113").Value = 3312
$ws.Range("L113").Value = 2665.75
$ws.Range("M113").Value = -1142
$ws.Range("N113").Value = -7005.75
$ws.Range("H122").Value = 2161.923
$ws.Range("I122").Value = 2161.923
$ws.Range("K122").Value = 6485.768999999999
$ws.Range("M122").Value = -4035.768999999999
$ws.Range("H126").Value = 18280
$ws.Range("J126").Value = 107007
$ws.Range("L126").Value = 321021
$ws.Range("N126").Value = -325961

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 599
$ws.Range("I23").Value = 573.75
$ws.Range("J23").Value = 700
$ws.Range("K23").Value = 573.75
$ws.Range("L23").Value = 700
$ws.Range("M23").Value = -343.75
$ws.Range("N23").Value = -1160
$ws.Range("H40").Value = 8084.6924
$ws.Range("J40").Value = 19974.25
$ws.Range("L40").Value = 19974.25
$ws.Range("N40").Value = -20246.25
$ws.Range("H68").Value = 3308.375
$ws.Range("I68").Value = 3266.5715
$ws.Range("J68").Value = 3601
$ws.Range("K68").Value = 3266.5715
$ws.Range("L68").Value = 3601
$ws.Range("M68").Value = -2517.5715
$ws.Range("N68").Value = -5099
$ws.Range("H71").Value = 3308.375
$ws.Range("I71").Value = 3266.5715
$ws.Range("J71").Value = 3601
$ws.Range("K71").Value = 16332.8575
$ws.Range("L71").Value = 18005
$ws.Range("M71").Value = -12588.8575
$ws.Range("N71").Value = -25493
$ws.Range("H75").Value = 53550.668
$ws.Range("I75").Value = 17826
$ws.Range("J75").Value = 125000
$ws.Range("K75").Value = 17826
$ws.Range("L75").Value = 125000
$ws.Range("M75").Value = -16890
$ws.Range("N75").Value = -126872
$ws.Range("H78").Value = 53550.668
$ws.Range("I78").Value = 17826
$ws.Range("J78").Value = 125000
$ws.Range("K78").Value = 53478
$ws.Range("L78").Value = 375000
$ws.Range("M78").Value = -48798
$ws.Range("N78").Value = -384360
$ws.Range("H81").Value = 112998.5
$ws.Range("J81").Value = 112998.5
$ws.Range("L81").Value = 112998.5
$ws.Range("N81").Value = -114994.5
$ws.Range("H84").Value = 112998.5
$ws.Range("J84").Value = 112998.5
$ws.Range("L84").Value = 338995.5
$ws.Range("N84").Value = -348979.5
$ws.Range("H97").Value = 34977.43
$ws.Range("J97").Value = 34977.43
$ws.Range("L97").Value = 34977.43
$ws.Range("N97").Value = -36959.43
$ws.Range("H122").Value = 11526.772
$ws.Range("I122").Value = 12645.77
$ws.Range("J122").Value = 9910.444
$ws.Range("K122").Value = 37937.31
$ws.Range("L122").Value = 29731.332
$ws.Range("M122").Value = -35487.31
$ws.Range("N122").Value = -34631.33199999999
$ws.Range("H132").Value = 3069.1777
$ws.Range("I132").Value = 1976.7241
$ws.Range("K132").Value = 5930.1723
$ws.Range("M132").Value = -3400.1723
$ws.Range("H136").Value = 3953.5894
$ws.Range("I136").Value = 3908.2654
$ws.Range("K136").Value = 11724.7962
$ws.Range("M136").Value = -9174.796200000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 30006
$ws.Range("I18").Value = 30006
$ws.Range("K18").Value = 30006
$ws.Range("M18").Value = -29833
$ws.Range("H54").Value = 30000
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H62").Value = 137026
$ws.Range("I62").Value = 137026
$ws.Range("K62").Value = 137026
$ws.Range("M62").Value = -136402
$ws.Range("H65").Value = 137026
$ws.Range("I65").Value = 137026
$ws.Range("K65").Value = 685130
$ws.Range("M65").Value = -682010
$ws.Range("H80").Value = 45000
$ws.Range("J80").Value = 45000
$ws.Range("L80").Value = 45000
$ws.Range("N80").Value = -46996
$ws.Range("H83").Value = 45000
$ws.Range("J83").Value = 45000
$ws.Range("L83").Value = 135000
$ws.Range("N83").Value = -144984
$ws.Range("H95").Value = 68333
$ws.Range("J95").Value = 68333
$ws.Range("L95").Value = 68333
$ws.Range("N95").Value = -73825
$ws.Range("H122").Value = 4155.5
$ws.Range("I122").Value = 4291.1
$ws.Range("K122").Value = 12873.3
$ws.Range("M122").Value = -10423.3
$ws.Range("H126").Value = 5988.2266
$ws.Range("I126").Value = 5321.7354
$ws.Range("K126").Value = 15965.2062
$ws.Range("M126").Value = -13495.2062
$ws.Range("H130").Value = 59750
$ws.Range("J130").Value = 59750
$ws.Range("L130").Value = 59750
$ws.Range("N130").Value = -69790
$ws.Range("H132").Value = 1641.1875
$ws.Range("I132").Value = 1578.6511
$ws.Range("K132").Value = 4735.9533
$ws.Range("M132").Value = -2205.9533
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

